$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + 1h volume/change %) per upstream refresh.
# D-column numeric-looking values need NumberFormat forced to Text first so Excel
# does not reinterpret the price string as a number (which would also drop
# significant trailing zeros, e.g. "163.70" -> 163.7).

$ws.Range('D2').Value = '72.943.77'
$ws.Range('E2').Value = '  +2.20%  '
$ws.Range('D3').Value = '3.996.26'
$ws.Range('E3').Value = '  +0.82%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '622.21'
$ws.Range('E5').Value = '  +15.66%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '163.70'
$ws.Range('E6').Value = '  +9.41%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.688'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.760'
$ws.Range('E9').Value = '  +1.94%  '
$ws.Range('E10').Value = '  +0.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '54.43'
$ws.Range('E11').Value = '  -1.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000320'
$ws.Range('E12').Value = '  -0.59%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.15'
$ws.Range('E13').Value = '  +3.90%  '
$ws.Range('D14').Value = '4.639.27'
$ws.Range('E14').Value = '  +1.07%  '
$ws.Range('D15').Value = '4.000.37'
$ws.Range('E15').Value = '  +1.14%  '
$ws.Range('E16').Value = '  +8.21%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.19'
$ws.Range('E17').Value = '  +1.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '20.72'
$ws.Range('E18').Value = '  +0.72%  '
$ws.Range('E19').Value = '  +0.41%  '
$ws.Range('D20').Value = '72.654.65'
$ws.Range('E20').Value = '  +1.97%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '441.71'
$ws.Range('E21').Value = '  +2.81%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.98'
$ws.Range('E22').Value = '  +17.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '97.16'
$ws.Range('E23').Value = '  -0.41%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.46'
$ws.Range('E24').Value = '  -3.72%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '14.48'
$ws.Range('E25').Value = '  -1.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.31'
$ws.Range('E26').Value = '  +8.36%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.37'
$ws.Range('E27').Value = '  -0.42%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.62'
$ws.Range('E28').Value = '  -1.58%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.96'
$ws.Range('E29').Value = '  +0.95%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.53'
$ws.Range('E30').Value = '  -0.67%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.74'
$ws.Range('E31').Value = '  -2.62%  '
$ws.Range('E32').Value = '  +4.64%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.132'
$ws.Range('E33').Value = '  +0.10%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '72.17'
$ws.Range('E34').Value = '  +9.82%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '48.31'
$ws.Range('E35').Value = '  -4.93%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '649.53'
$ws.Range('E36').Value = '  -4.81%  '
$ws.Range('D37').Value = '0.0₃0911'
$ws.Range('E37').Value = '  +10.54%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.440'
$ws.Range('E38').Value = '  -0.52%  '
$ws.Range('E39').Value = '  -0.40%  '
$ws.Range('E40').Value = '  -0.55%  '
$ws.Range('E41').Value = '  -0.18%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.36'
$ws.Range('E42').Value = '  +4.59%  '
$ws.Range('E43').Value = '  +0.22%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0492'
$ws.Range('E44').Value = '  +1.46%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.67'
$ws.Range('E45').Value = '  +3.26%  '
$ws.Range('E46').Value = '  +1.05%  '
$ws.Range('E47').Value = '  -0.16%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.42'
$ws.Range('E48').Value = '  +1.69%  '
$ws.Range('D49').Value = '2.929.52'
$ws.Range('E49').Value = '  +12.03%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.09'
$ws.Range('E50').Value = '  +2.67%  '
$ws.Range('E51').Value = '  +4.40%  '
